$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178617477416992
$ws.Range("B1").Value = 2.186071157455444
$ws.Range("C1").Value = 3.261564493179321
$ws.Range("D1").Value = 3.616487741470337
$ws.Range("E1").Value = 1.162784934043884
